# Refresh the cryptos list (coin rank shuffles + updated price/volume figures).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin / link / percentage cells -- plain text, Excel will not mis-parse these as numbers.
$textCells = @{
    "D2" = "69.472.74"
    "E2" = "  +2.08%  "
    "D3" = "3.390.63"
    "E3" = "  +1.65%  "
    "E4" = "  +0.21%  "
    "E5" = "  +0.17%  "
    "E6" = "  +2.42%  "
    "B7" = "USDC"
    "C7" = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
    "E7" = "  +0.14%  "
    "B8" = "XRP"
    "C8" = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
    "E8" = "  +0.51%  "
    "E9" = "  +9.84%  "
    "E10" = "  +1.44%  "
    "E11" = "  +1.22%  "
    "E12" = "  +4.44%  "
    "E13" = "  -2.58%  "
    "E14" = "  +2.94%  "
    "D15" = "3.952.36"
    "E15" = "  +1.99%  "
    "D16" = "69.463.27"
    "E16" = "  +2.05%  "
    "D17" = "3.405.41"
    "E17" = "  +2.71%  "
    "E18" = "  +1.37%  "
    "E19" = "  +1.32%  "
    "E20" = "  +1.55%  "
    "E21" = "  +1.51%  "
    "E22" = "  +1.66%  "
    "E23" = "  -0.71%  "
    "E24" = "  +1.65%  "
    "E25" = "  +0.30%  "
    "E26" = "  +0.92%  "
    "E27" = "  +3.07%  "
    "E28" = "  +1.92%  "
    "E29" = "  +3.14%  "
    "E30" = "  -0.85%  "
    "E31" = "  +0.96%  "
    "E32" = "  -2.74%  "
    "B33" = "dogwifhat"
    "C33" = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
    "E33" = "  +8.19%  "
    "B34" = "Hedera"
    "C34" = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
    "E34" = "  +1.19%  "
    "E35" = "  +2.27%  "
    "E36" = "  +0.12%  "
    "D37" = "3.662.85"
    "E37" = "  -2.44%  "
    "E38" = "  +4.27%  "
    "E39" = "  +0.52%  "
    "D40" = "0.0₃0728"
    "E40" = "  +7.57%  "
    "E41" = "  +4.51%  "
    "E42" = "  +3.49%  "
    "B43" = "ApeXProtocol"
    "C43" = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
    "E43" = "  +2.33%  "
    "B44" = "VeChain"
    "C44" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "E44" = "  +5.03%  "
    "B45" = "TheGraph"
    "C45" = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
    "E45" = "  +1.26%  "
    "E46" = "  +1.36%  "
    "E47" = "  +0.53%  "
    "B48" = "FirstDigitalUSD"
    "C48" = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
    "E48" = "  +0.04%  "
    "B49" = "Mantle"
    "C49" = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
    "E49" = "  +4.44%  "
    "E50" = "  -0.73%  "
    "E51" = "  -1.68%  "
}
foreach ($ref in $textCells.Keys) {
    $ws.Range($ref).Value = $textCells[$ref]
}

# Price cells whose text LOOKS numeric (e.g. "1.00", "0.200") -- must stay literal text so
# trailing zeros / exact digits survive instead of Excel auto-converting to a number.
$numericTextCells = @{
    "D5" = "584.56"
    "D6" = "180.18"
    "D7" = "1.00"
    "D8" = "0.596"
    "D9" = "0.200"
    "D10" = "0.591"
    "D11" = "48.28"
    "D12" = "0.0000286"
    "D13" = "683.09"
    "D14" = "8.65"
    "D19" = "17.71"
    "D20" = "11.32"
    "D21" = "0.910"
    "D22" = "17.26"
    "D23" = "5.35"
    "D24" = "102.49"
    "D25" = "3.92"
    "D26" = "2.71"
    "D27" = "9.71"
    "D28" = "33.87"
    "D29" = "8.80"
    "D30" = "6.94"
    "D31" = "11.11"
    "D32" = "553.97"
    "D33" = "3.61"
    "D34" = "0.107"
    "D35" = "58.35"
    "D38" = "0.140"
    "D39" = "35.66"
    "D41" = "3.29"
    "D42" = "2.71"
    "D43" = "3.37"
    "D44" = "0.0429"
    "D45" = "0.338"
    "D46" = "2.67"
    "D48" = "1.00"
    "D49" = "1.38"
    "D50" = "129.23"
    "D51" = "2.67"
}
foreach ($ref in $numericTextCells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"        # force Text format so the assignment is not parsed as a number
    $cell.Value = $numericTextCells[$ref]
    $cell.ClearFormats()            # restore original (default) cell style/format
}
